$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 7400
$ws.Range("J54").Value = 7400
$ws.Range("L54").Value = 7400
$ws.Range("N54").Value = -8372

$ws.Range("H62").Value = 6189756
$ws.Range("I62").Value = 9268800
$ws.Range("J62").Value = 31666.666
$ws.Range("K62").Value = 9268800
$ws.Range("L62").Value = 31666.666
$ws.Range("M62").Value = -9268176
$ws.Range("N62").Value = -32914.666

$ws.Range("H65").Value = 6189756
$ws.Range("I65").Value = 9268800
$ws.Range("J65").Value = 31666.666
$ws.Range("K65").Value = 46344000
$ws.Range("L65").Value = 158333.33
$ws.Range("M65").Value = -46340880
$ws.Range("N65").Value = -164573.33

$ws.Range("H133").Value = 49935
$ws.Range("J133").Value = 49935
$ws.Range("L133").Value = 49935
$ws.Range("N133").Value = -60055


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 52415.266
$ws.Range("I32").Value = 10097.76
$ws.Range("J32").Value = 264002.8
$ws.Range("K32").Value = 10097.76
$ws.Range("L32").Value = 264002.8
$ws.Range("M32").Value = -9810.76
$ws.Range("N32").Value = -264576.8

$ws.Range("H53").Value = 3000
$ws.Range("I53").Value = 3000
$ws.Range("K53").Value = 3000
$ws.Range("M53").Value = -2318

$ws.Range("H74").Value = 4597.9736
$ws.Range("I74").Value = 1024.6522
$ws.Range("J74").Value = 10077.066
$ws.Range("K74").Value = 1024.6522
$ws.Range("L74").Value = 10077.066
$ws.Range("M74").Value = -150.6522
$ws.Range("N74").Value = -11825.066

$ws.Range("H77").Value = 4597.9736
$ws.Range("I77").Value = 1024.6522
$ws.Range("J77").Value = 10077.066
$ws.Range("K77").Value = 5123.261
$ws.Range("L77").Value = 50385.33
$ws.Range("M77").Value = -755.2610000000004
$ws.Range("N77").Value = -59121.33

$ws.Range("H122").Value = 6487
$ws.Range("I122").Value = 6487
$ws.Range("K122").Value = 19461
$ws.Range("M122").Value = -17011

$ws.Range("H133").Value = 53200
$ws.Range("J133").Value = 53200
$ws.Range("L133").Value = 53200
$ws.Range("N133").Value = -58260

$ws.Range("H139").Value = 48735
$ws.Range("J139").Value = 48735
$ws.Range("L139").Value = 48735
$ws.Range("N139").Value = -59015


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 59450
$ws.Range("J59").Value = 59450
$ws.Range("L59").Value = 59450
$ws.Range("N59").Value = -61144


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 12824067
$ws.Range("I99").Value = 3614
$ws.Range("J99").Value = 55558908
$ws.Range("K99").Value = 3614
$ws.Range("L99").Value = 55558908
$ws.Range("M99").Value = -2116
$ws.Range("N99").Value = -55561904

$ws.Range("H105").Value = 835.5484
$ws.Range("I105").Value = 683.4167
$ws.Range("K105").Value = 683.4167
$ws.Range("M105").Value = 1063.5833

$ws.Range("H126").Value = 12824067
$ws.Range("I126").Value = 3614
$ws.Range("J126").Value = 55558908
$ws.Range("K126").Value = 10842
$ws.Range("L126").Value = 166676724
$ws.Range("M126").Value = -8372
$ws.Range("N126").Value = -166681664

$ws.Range("H132").Value = 4168842.2
$ws.Range("I132").Value = 5556912
$ws.Range("K132").Value = 16670736
$ws.Range("M132").Value = -16668206

$ws.Range("H140").Value = 79826.664
$ws.Range("J140").Value = 79826.664
$ws.Range("L140").Value = 79826.664
$ws.Range("N140").Value = -90186.664

$ws.Range("H141").Value = 221429.89
$ws.Range("J141").Value = 224982.05
$ws.Range("L141").Value = 224982.05
$ws.Range("N141").Value = -235342.05


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H28").Value = 1049.3334
$ws.Range("I28").Value = 443.33334
$ws.Range("J28").Value = 1655.3334
$ws.Range("K28").Value = 1330.00002
$ws.Range("L28").Value = 4966.0002
$ws.Range("M28").Value = -1098.00002
$ws.Range("N28").Value = -5430.0002

$ws.Range("H75").Value = 3000
$ws.Range("J75").Value = 3000
$ws.Range("L75").Value = 9000
$ws.Range("N75").Value = -10996

$ws.Range("H78").Value = 3000
$ws.Range("J78").Value = 3000
$ws.Range("L78").Value = 27000
$ws.Range("N78").Value = -36984


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 6000
$ws.Range("I52").Value = 3000
$ws.Range("K52").Value = 3000
$ws.Range("M52").Value = -2741

$ws.Range("H55").Value = 4500
$ws.Range("I55").Value = 3000
$ws.Range("J55").Value = 6000
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 6000
$ws.Range("M55").Value = -2673
$ws.Range("N55").Value = -6654

$ws.Range("H102").Value = 6145.75
$ws.Range("I102").Value = 3957.4546
$ws.Range("J102").Value = 10960
$ws.Range("K102").Value = 3957.4546
$ws.Range("L102").Value = 10960
$ws.Range("M102").Value = -2335.4546
$ws.Range("N102").Value = -14204

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 170
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 2627.3076
$ws.Range("I122").Value = 2700.4167
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 8101.250100000001
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -5651.250100000001
$ws.Range("N122").Value = -10150

$ws.Range("H138").Value = 89500
$ws.Range("J138").Value = 89500
$ws.Range("L138").Value = 89500
$ws.Range("N138").Value = -99780

$ws.Range("H139").Value = 56363
$ws.Range("J139").Value = 56363
$ws.Range("L139").Value = 56363
$ws.Range("N139").Value = -66643


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 49500
$ws.Range("I39").Value = 49500
$ws.Range("K39").Value = 49500
$ws.Range("M39").Value = -49040

$ws.Range("H40").Value = 2984.4211
$ws.Range("I40").Value = 1867.3334
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 1867.3334
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = -1731.3334
$ws.Range("N40").Value = -3772

$ws.Range("H48").Value = 10760

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H136").Value = 14071
$ws.Range("J136").Value = 12237.333
$ws.Range("L136").Value = 36711.999
$ws.Range("N136").Value = -41811.999


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4812.875
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 4812.875
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 9625.75
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -11747.75

$ws.Range("H84").Value = 4812.875
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 4812.875
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 48128.75
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -58736.75

$ws.Range("H126").Value = 100744.3
$ws.Range("I126").Value = 143234.72
$ws.Range("J126").Value = 1600
$ws.Range("K126").Value = 429704.16
$ws.Range("L126").Value = 4800
$ws.Range("M126").Value = -427234.16
$ws.Range("N126").Value = -9740

$ws.Range("H132").Value = 3966.2104
$ws.Range("I132").Value = 3746.0386
$ws.Range("J132").Value = 4443.25
$ws.Range("K132").Value = 11238.1158
$ws.Range("L132").Value = 13329.75
$ws.Range("M132").Value = -8708.1158
$ws.Range("N132").Value = -18389.75

